# Tasks List report update:
#  - add a "Label" column (F) with a tag/category per task
#  - add two new tasks ("Tarea Simulación", "Tarea SO")
#  - mark several existing tasks as Done and renumber the "No Task" column
#  - widen columns B and F to fit the new content
#
# Due-Date values ("YYYY-MM-DD") must stay plain text, like the source
# workbook already stores them - so each date cell is pre-formatted as
# Text ("@") before the value is assigned (otherwise Excel's normal
# autodetection would silently convert it into a date serial number),
# and the Text number-format is cleared back off afterwards so the cell
# keeps the worksheet's default look.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# New "Label" header
$ws.Range("F2").Value = "Label"

# Row 3: Tarea Cultura -> now Done, labelled Personal
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = "Tarea Cultura"
$ws.Range("C3").Value = "Hacer el reporte de cultura"
$ws.Range("D3").Value = $true
Set-TextValue $ws.Range("E3") "2023-10-31"
$ws.Range("F3").Value = "Personal"

# Row 4: Tarea Simulación (new task)
$ws.Range("A4").Value = 6
$ws.Range("B4").Value = "Tarea Simulación"
$ws.Range("C4").Value = "Hacer las distribuciones"
$ws.Range("D4").Value = $false
Set-TextValue $ws.Range("E4") "2023-11-14"
$ws.Range("F4").Value = "Simulación, ITC"

# Row 5: Tarea Taller -> now Done, labelled Personal
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "Tarea Taller"
$ws.Range("C5").Value = "Hacer el Laboratorio de TBD"
$ws.Range("D5").Value = $true
Set-TextValue $ws.Range("E5") "2023-10-31"
$ws.Range("F5").Value = "Personal"

# Row 6: Tarea Topicos -> now Done, labelled Personal,ITC
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = "Tarea Topicos"
$ws.Range("C6").Value = "Hacer la tarea todo de topicos"
$ws.Range("D6").Value = $true
Set-TextValue $ws.Range("E6") "2023-11-17"
$ws.Range("F6").Value = "Personal,ITC"

# Row 7: Tarea SO (new task)
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "Tarea SO"
$ws.Range("C7").Value = "Montar maquina virtual"
$ws.Range("D7").Value = $true
Set-TextValue $ws.Range("E7") "2023-11-13"
$ws.Range("F7").Value = "ITC,Sistemas operativos"

# Widen columns B (longer titles) and F (new Label column) to fit content
$ws.Columns.Item(2).ColumnWidth = 13.59765625
$ws.Columns.Item(6).ColumnWidth = 18.54296875
